$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.52848219871521
$ws.Range("B1").Value = 1.694281578063965
$ws.Range("C1").Value = 1.81993567943573
$ws.Range("D1").Value = 1.471868872642517
$ws.Range("E1").Value = 1.27924108505249
